# Applies the "EC" update: removes/replaces the previous single-period
# account-statement row with three period rows (2505, 2504, 2503),
# updates the totals (Valor Mora total, Cant. Periodos) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new rows right below the existing data row (row 16) so the
#    signature block (previously rows 21/22) moves down to rows 23/24,
#    exactly mirroring what Excel does when a user inserts rows in the
#    middle of the sheet.
$ws.Rows("17:18").Insert()

# 2) Duplicate the formatting of the existing data row (row 16) into the
#    two freshly inserted rows, the same way a user would by copying row
#    16 and pasting it onto rows 17 and 18.
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))
$ws.Range("B16:J16").Copy($ws.Range("B18:J18"))

# 3) Set the "Periodo Mora" value for each of the three data rows.
$ws.Range("E16").Value = "2505"
$ws.Range("E17").Value = "2504"
$ws.Range("E18").Value = "2503"

# 4) Update the summary fields: 3 periods now reported, and the total
#    "Valor Mora" is the sum of the three periods' value (376144 * 3).
$ws.Range("F13").Value = 3
$ws.Range("E11").Value = $ws.Range("F16").Value() + $ws.Range("F17").Value() + $ws.Range("F18").Value()
